$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.834.23'
$ws.Range("E2").Value = '  +4.56%  '
$ws.Range("D3").Value = '3.064.56'
$ws.Range("E3").Value = '  +2.32%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'576.65"
$ws.Range("E5").Value = '  +2.58%  '
$ws.Range("D6").Value = "'142.87"
$ws.Range("E6").Value = '  +3.13%  '
$ws.Range("D8").Value = '3.055.00'
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("E10").Value = '  +4.85%  '
$ws.Range("D11").Value = "'5.47"
$ws.Range("E11").Value = '  +11.19%  '
$ws.Range("E12").Value = '  +1.84%  '
$ws.Range("E13").Value = '  +4.48%  '
$ws.Range("D14").Value = "'34.85"
$ws.Range("E14").Value = '  +3.17%  '
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").Value = '3.570.46'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = "'7.21"
$ws.Range("E17").Value = '  +3.01%  '
$ws.Range("D18").Value = '3.067.77'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = '61.843.18'
$ws.Range("E19").Value = '  +4.56%  '
$ws.Range("D20").Value = "'449.87"
$ws.Range("E20").Value = '  +6.01%  '
$ws.Range("D21").Value = "'13.94"
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("D22").Value = "'0.731"
$ws.Range("E22").Value = '  +2.69%  '
$ws.Range("D23").Value = "'7.32"
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").Value = "'81.83"
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +4.45%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  +3.84%  '
$ws.Range("D30").Value = "'8.07"
$ws.Range("E30").Value = '  +3.54%  '
$ws.Range("D31").Value = "'6.64"
$ws.Range("E31").Value = '  +8.23%  '
$ws.Range("D32").Value = "'26.60"
$ws.Range("E32").Value = '  +3.83%  '
$ws.Range("E33").Value = '  +8.04%  '
$ws.Range("D34").Value = '0.0₃0806'
$ws.Range("E34").Value = '  +5.02%  '
$ws.Range("E35").Value = '  +2.15%  '
$ws.Range("E36").Value = '  +5.41%  '
$ws.Range("E37").Value = '  +4.90%  '
$ws.Range("D38").Value = "'50.02"
$ws.Range("E38").Value = '  +2.08%  '
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = '  +7.03%  '
$ws.Range("D40").Value = "'8.80"
$ws.Range("E40").Value = '  +1.56%  '
$ws.Range("D41").Value = "'418.09"
$ws.Range("E41").Value = '  +3.93%  '
$ws.Range("E42").Value = '  +5.24%  '
$ws.Range("D43").Value = '2.771.96'
$ws.Range("E43").Value = '  +0.94%  '
$ws.Range("D44").Value = "'0.108"
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("E45").Value = '  +7.98%  '
$ws.Range("D46").Value = "'36.37"
$ws.Range("E46").Value = '  +11.42%  '
$ws.Range("E47").Value = '  +3.62%  '
$ws.Range("D49").Value = "'123.25"
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("D51").Value = "'24.08"
$ws.Range("E51").Value = '  +2.61%  '

# Reset style on cells that needed a quote-prefix so no stray number-format/quotePrefix
# style gets attached to the cell (matches original plain default styling).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
